$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 40.497736
$ws.Range("H2").Value = 80.99547200000001
$ws.Range("I2").Value = 0.4338277436889532
$ws.Range("J2").Value = 0.3549025979779166
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 2.675209
$ws.Range("N2").Value = 5.350417999999999
$ws.Range("O2").Value = 0.6038885578943551
$ws.Range("P2").Value = 0.5040571546736343
$ws.Range("Q2").Value = 108.339907826824
$ws.Range("R2").Value = 433.359631307296
$ws.Range("S2").Value = 0.2619836105108839
$ws.Range("T2").Value = 0.1788911937230294

$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 40.497736
$ws.Range("H3").Value = 80.99547200000001
$ws.Range("I3").Value = 0.4338277436889532
$ws.Range("J3").Value = 0.3549025979779166
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.754762333333333
$ws.Range("N3").Value = 5.264286999999999
$ws.Range("O3").Value = 0.3961114421056449
$ws.Range("P3").Value = 0.4959428453263656
$ws.Range("Q3").Value = 71.06390171807733
$ws.Range("R3").Value = 426.383410308464
$ws.Range("S3").Value = 0.1718441331780693
$ws.Range("T3").Value = 0.1760114042548872

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.464838333333333
$ws.Range("H4").Value = 19.394515
$ws.Range("I4").Value = 0.06925390168635236
$ws.Range("J4").Value = 0.08498208097388052
$ws.Range("K4").Value = 2
$ws.Range("M4").Value = 2.675209
$ws.Range("N4").Value = 5.350417999999999
$ws.Range("O4").Value = 0.6038885578943551
$ws.Range("P4").Value = 0.5040571546736343
$ws.Range("Q4").Value = 17.29479369287833
$ws.Range("R4").Value = 103.76876215727
$ws.Range("S4").Value = 0.04182163881792877
$ws.Range("T4").Value = 0.04283582593393861

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.464838333333333
$ws.Range("H5").Value = 19.394515
$ws.Range("I5").Value = 0.06925390168635236
$ws.Range("J5").Value = 0.08498208097388052
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.754762333333333
$ws.Range("N5").Value = 5.264286999999999
$ws.Range("O5").Value = 0.3961114421056449
$ws.Range("P5").Value = 0.4959428453263656
$ws.Range("Q5").Value = 11.34425479842277
$ws.Range("R5").Value = 102.098293185805
$ws.Range("S5").Value = 0.02743226286842359
$ws.Range("T5").Value = 0.0421462550399419

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.7642823333333334
$ws.Range("H6").Value = 2.292847
$ws.Range("I6").Value = 0.008187294228282994
$ws.Range("J6").Value = 0.01004670183372562
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 2.675209
$ws.Range("N6").Value = 5.350417999999999
$ws.Range("O6").Value = 0.6038885578943551
$ws.Range("P6").Value = 0.5040571546736343
$ws.Range("Q6").Value = 2.044614976674333
$ws.Range("R6").Value = 12.267689860046
$ws.Range("S6").Value = 0.004944213304574595
$ws.Range("T6").Value = 0.005064111940162121

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.7642823333333334
$ws.Range("H7").Value = 2.292847
$ws.Range("I7").Value = 0.008187294228282994
$ws.Range("J7").Value = 0.01004670183372562
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.754762333333333
$ws.Range("N7").Value = 5.264286999999999
$ws.Range("O7").Value = 0.3961114421056449
$ws.Range("P7").Value = 0.4959428453263656
$ws.Range("Q7").Value = 1.341133850565444
$ws.Range("R7").Value = 12.070204655089
$ws.Range("S7").Value = 0.0032430809237084
$ws.Range("T7").Value = 0.004982589893563499

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.535784
$ws.Range("H8").Value = 1.607352
$ws.Range("I8").Value = 0.00573952983012784
$ws.Range("J8").Value = 0.007043028290087626
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 2.675209
$ws.Range("N8").Value = 5.350417999999999
$ws.Range("O8").Value = 0.6038885578943551
$ws.Range("P8").Value = 0.5040571546736343
$ws.Range("Q8").Value = 1.433334178856
$ws.Range("R8").Value = 8.600005073136
$ws.Range("S8").Value = 0.003466036392107534
$ws.Range("T8").Value = 0.003550088800187481

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.535784
$ws.Range("H9").Value = 1.607352
$ws.Range("I9").Value = 0.00573952983012784
$ws.Range("J9").Value = 0.007043028290087626
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.754762333333333
$ws.Range("N9").Value = 5.264286999999999
$ws.Range("O9").Value = 0.3961114421056449
$ws.Range("P9").Value = 0.4959428453263656
$ws.Range("Q9").Value = 0.9401735820026667
$ws.Range("R9").Value = 8.461562238023999
$ws.Range("S9").Value = 0.002273493438020306
$ws.Range("T9").Value = 0.003492939489900145

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 33.754358
$ws.Range("H10").Value = 101.263074
$ws.Range("I10").Value = 0.3615900150766247
$ws.Range("J10").Value = 0.4437103353361533
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 2.675209
$ws.Range("N10").Value = 5.350417999999999
$ws.Range("O10").Value = 0.6038885578943551
$ws.Range("P10").Value = 0.5040571546736343
$ws.Range("Q10").Value = 90.299962310822
$ws.Range("R10").Value = 541.7997738649319
$ws.Range("S10").Value = 0.218360072753621
$ws.Range("T10").Value = 0.2236553691288256

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 33.754358
$ws.Range("H11").Value = 101.263074
$ws.Range("I11").Value = 0.3615900150766247
$ws.Range("J11").Value = 0.4437103353361533
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.754762333333333
$ws.Range("N11").Value = 5.264286999999999
$ws.Range("O11").Value = 0.3961114421056449
$ws.Range("P11").Value = 0.4959428453263656
$ws.Range("Q11").Value = 59.23087600424866
$ws.Range("R11").Value = 533.0778840382379
$ws.Range("S11").Value = 0.1432299423230037
$ws.Range("T11").Value = 0.2200549662073277

$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 11.332808
$ws.Range("H12").Value = 22.665616
$ws.Range("I12").Value = 0.1214015154896589
$ws.Range("J12").Value = 0.09931525558823626
$ws.Range("K12").Value = 2
$ws.Range("M12").Value = 2.675209
$ws.Range("N12").Value = 5.350417999999999
$ws.Range("O12").Value = 0.6038885578943551
$ws.Range("P12").Value = 0.5040571546736343
$ws.Range("Q12").Value = 30.317629956872
$ws.Range("R12").Value = 121.270519827488
$ws.Range("S12").Value = 0.07331298611523936
$ws.Range("T12").Value = 0.05006056514749112

$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 11.332808
$ws.Range("H13").Value = 22.665616
$ws.Range("I13").Value = 0.1214015154896589
$ws.Range("J13").Value = 0.09931525558823626
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.754762333333333
$ws.Range("N13").Value = 5.264286999999999
$ws.Range("O13").Value = 0.3961114421056449
$ws.Range("P13").Value = 0.4959428453263656
$ws.Range("Q13").Value = 19.88638460929866
$ws.Range("R13").Value = 119.318307655792
$ws.Range("S13").Value = 0.04808852937441959
$ws.Range("T13").Value = 0.04925469044074512
